$wb = $excel.ActiveWorkbook

# Update CAPEX value on the "Inputs" sheet (cell C2) from 100,000,000 to 400,000,000.
# All other changed cells in the diff are formulas that depend on this value
# (Simluation_Tool_Example!B2, E2, B23, B24 and Outputs!C3, C4) and will
# recalculate automatically.
$inputs = $wb.Worksheets.Item("Inputs")
$inputs.Range("C2").Value = 400000000
